{"js": "// The \"1TI\" book-intro heading (Heading 2) is immediately followed by a\n// short paragraph whose only visible content is the italic run \"1 Timothy\"\n// (the English book title, shown under the localized \"1TI\" abbreviation).\n// This paragraph is removed entirely, per the resource update.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text,items/style,items/font/italic\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.style === \"Normal\" && p.text.trim() === \"1 Timothy\" && p.font.italic) {\n    p.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# The \"1TI\" book-intro heading (Heading 2) is immediately followed by a\n# short paragraph whose only visible content is the italic run \"1 Timothy\"\n# (the English book title, shown under the localized \"1TI\" abbreviation).\n# This paragraph is removed entirely, per the resource update.\n$d = $word.ActiveDocument\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $r = $p.Range\n    $text = $r.Text.Trim()\n    $styleName = $r.Style.NameLocal\n    $isItalic = $r.Font.Italic\n\n    if ($styleName -eq \"Normal\" -and $text -eq \"1 Timothy\" -and $isItalic -eq -1) {\n        $r.Delete()\n    }\n}\n"}
